$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69
$ws.Range('A69').NumberFormat = '@'
$ws.Range('A69').Value = '7180'
$ws.Range('A69').Style = 'Normal'
$ws.Range('B69').NumberFormat = '@'
$ws.Range('B69').Value = '9/8/2025'
$ws.Range('B69').Style = 'Normal'
$ws.Range('C69').NumberFormat = '@'
$ws.Range('C69').Value = 'GORRITI 4417'
$ws.Range('C69').Style = 'Normal'
$ws.Range('D69').NumberFormat = '@'
$ws.Range('D69').Value = '14'
$ws.Range('D69').Style = 'Normal'
$ws.Range('E69').NumberFormat = '@'
$ws.Range('E69').Value = '809526157'
$ws.Range('E69').Style = 'Normal'
$ws.Range('F69').NumberFormat = '@'
$ws.Range('F69').Value = 'PEBCOM'
$ws.Range('F69').Style = 'Normal'
$ws.Range('G69').NumberFormat = '@'
$ws.Range('G69').Value = 'Pendiente'
$ws.Range('G69').Style = 'Normal'
$ws.Range('H69').NumberFormat = '@'
$ws.Range('H69').Value = 'Picada'
$ws.Range('H69').Style = 'Normal'
$ws.Range('J69').NumberFormat = '@'
$ws.Range('J69').Value = 'Cambio'
$ws.Range('J69').Style = 'Normal'
$ws.Range('K69').NumberFormat = '@'
$ws.Range('K69').Value = 'Sin equipos'
$ws.Range('K69').Style = 'Normal'
$ws.Range('L69').NumberFormat = '@'
$ws.Range('L69').Value = 'Pasante'
$ws.Range('L69').Style = 'Normal'
$ws.Range('O69').NumberFormat = '@'
$ws.Range('O69').Value = 'Palermo'
$ws.Range('O69').Style = 'Normal'
$ws.Range('P69').NumberFormat = '@'
$ws.Range('P69').Value = 'Capital Sur'
$ws.Range('P69').Style = 'Normal'
$ws.Range('I69').Value = 1
$ws.Range('M69').Value = -58.425358
$ws.Range('N69').Value = -34.593308

# Row 70
$ws.Range('A70').NumberFormat = '@'
$ws.Range('A70').Value = '7186'
$ws.Range('A70').Style = 'Normal'
$ws.Range('B70').NumberFormat = '@'
$ws.Range('B70').Value = '9/8/2025'
$ws.Range('B70').Style = 'Normal'
$ws.Range('C70').NumberFormat = '@'
$ws.Range('C70').Value = 'NICARAGUA 5510'
$ws.Range('C70').Style = 'Normal'
$ws.Range('D70').NumberFormat = '@'
$ws.Range('D70').Value = '14'
$ws.Range('D70').Style = 'Normal'
$ws.Range('E70').NumberFormat = '@'
$ws.Range('E70').Value = '809526162'
$ws.Range('E70').Style = 'Normal'
$ws.Range('F70').NumberFormat = '@'
$ws.Range('F70').Value = 'PEBCOM'
$ws.Range('F70').Style = 'Normal'
$ws.Range('G70').NumberFormat = '@'
$ws.Range('G70').Value = 'Pendiente'
$ws.Range('G70').Style = 'Normal'
$ws.Range('H70').NumberFormat = '@'
$ws.Range('H70').Value = 'Picada'
$ws.Range('H70').Style = 'Normal'
$ws.Range('J70').NumberFormat = '@'
$ws.Range('J70').Value = 'Cambio'
$ws.Range('J70').Style = 'Normal'
$ws.Range('K70').NumberFormat = '@'
$ws.Range('K70').Value = 'Sin equipos'
$ws.Range('K70').Style = 'Normal'
$ws.Range('L70').NumberFormat = '@'
$ws.Range('L70').Value = 'Pasante'
$ws.Range('L70').Style = 'Normal'
$ws.Range('O70').NumberFormat = '@'
$ws.Range('O70').Value = 'Palermo'
$ws.Range('O70').Style = 'Normal'
$ws.Range('P70').NumberFormat = '@'
$ws.Range('P70').Value = 'Capital Sur'
$ws.Range('P70').Style = 'Normal'
$ws.Range('I70').Value = 1
$ws.Range('M70').Value = -58.432726
$ws.Range('N70').Value = -34.582328

# Row 71
$ws.Range('A71').NumberFormat = '@'
$ws.Range('A71').Value = '7194'
$ws.Range('A71').Style = 'Normal'
$ws.Range('B71').NumberFormat = '@'
$ws.Range('B71').Value = '9/8/2025'
$ws.Range('B71').Style = 'Normal'
$ws.Range('C71').NumberFormat = '@'
$ws.Range('C71').Value = 'CASEROS AV. 2032'
$ws.Range('C71').Style = 'Normal'
$ws.Range('D71').NumberFormat = '@'
$ws.Range('D71').Value = '4'
$ws.Range('D71').Style = 'Normal'
$ws.Range('E71').NumberFormat = '@'
$ws.Range('E71').Value = 'ICD30709299'
$ws.Range('E71').Style = 'Normal'
$ws.Range('F71').NumberFormat = '@'
$ws.Range('F71').Value = 'PEBCOM'
$ws.Range('F71').Style = 'Normal'
$ws.Range('G71').NumberFormat = '@'
$ws.Range('G71').Value = 'Pendiente'
$ws.Range('G71').Style = 'Normal'
$ws.Range('H71').NumberFormat = '@'
$ws.Range('H71').Value = 'Aplomo '
$ws.Range('H71').Style = 'Normal'
$ws.Range('J71').NumberFormat = '@'
$ws.Range('J71').Value = 'Aplomo'
$ws.Range('J71').Style = 'Normal'
$ws.Range('K71').NumberFormat = '@'
$ws.Range('K71').Value = 'Sin equipos'
$ws.Range('K71').Style = 'Normal'
$ws.Range('L71').NumberFormat = '@'
$ws.Range('L71').Value = 'Terminal'
$ws.Range('L71').Style = 'Normal'
$ws.Range('O71').NumberFormat = '@'
$ws.Range('O71').Value = 'San Telmo'
$ws.Range('O71').Style = 'Normal'
$ws.Range('P71').NumberFormat = '@'
$ws.Range('P71').Value = 'Capital Sur'
$ws.Range('P71').Style = 'Normal'
$ws.Range('I71').Value = 1
$ws.Range('M71').Value = -58.390906
$ws.Range('N71').Value = -34.634312

# Row 72
$ws.Range('A72').NumberFormat = '@'
$ws.Range('A72').Value = '-586'
$ws.Range('A72').Style = 'Normal'
$ws.Range('B72').NumberFormat = '@'
$ws.Range('B72').Value = '9/8/2025'
$ws.Range('B72').Style = 'Normal'
$ws.Range('C72').NumberFormat = '@'
$ws.Range('C72').Value = 'Franklin 666'
$ws.Range('C72').Style = 'Normal'
$ws.Range('D72').NumberFormat = '@'
$ws.Range('D72').Value = '6'
$ws.Range('D72').Style = 'Normal'
$ws.Range('E72').NumberFormat = '@'
$ws.Range('E72').Value = 'ICD30709119'
$ws.Range('E72').Style = 'Normal'
$ws.Range('F72').NumberFormat = '@'
$ws.Range('F72').Value = 'PEBCOM'
$ws.Range('F72').Style = 'Normal'
$ws.Range('G72').NumberFormat = '@'
$ws.Range('G72').Value = 'Pendiente'
$ws.Range('G72').Style = 'Normal'
$ws.Range('H72').NumberFormat = '@'
$ws.Range('H72').Value = 'Traspasar redes a la columna de telecentro y desmontar la picada'
$ws.Range('H72').Style = 'Normal'
$ws.Range('J72').NumberFormat = '@'
$ws.Range('J72').Value = 'Desmonte'
$ws.Range('J72').Style = 'Normal'
$ws.Range('K72').NumberFormat = '@'
$ws.Range('K72').Value = 'Sin equipos'
$ws.Range('K72').Style = 'Normal'
$ws.Range('L72').NumberFormat = '@'
$ws.Range('L72').Value = 'Pasante'
$ws.Range('L72').Style = 'Normal'
$ws.Range('O72').NumberFormat = '@'
$ws.Range('O72').Value = 'Almagro'
$ws.Range('O72').Style = 'Normal'
$ws.Range('P72').NumberFormat = '@'
$ws.Range('P72').Value = 'Capital Sur'
$ws.Range('P72').Style = 'Normal'
$ws.Range('I72').Value = 1
$ws.Range('M72').Value = -58.441362
$ws.Range('N72').Value = -34.607784

# Row 73
$ws.Range('A73').NumberFormat = '@'
$ws.Range('A73').Value = '7224'
$ws.Range('A73').Style = 'Normal'
$ws.Range('B73').NumberFormat = '@'
$ws.Range('B73').Value = '9/16/2025'
$ws.Range('B73').Style = 'Normal'
$ws.Range('C73').NumberFormat = '@'
$ws.Range('C73').Value = 'CABILDO AV. 3950'
$ws.Range('C73').Style = 'Normal'
$ws.Range('D73').NumberFormat = '@'
$ws.Range('D73').Value = '12'
$ws.Range('D73').Style = 'Normal'
$ws.Range('E73').NumberFormat = '@'
$ws.Range('E73').Value = '809784515'
$ws.Range('E73').Style = 'Normal'
$ws.Range('F73').NumberFormat = '@'
$ws.Range('F73').Value = 'PEBCOM'
$ws.Range('F73').Style = 'Normal'
$ws.Range('G73').NumberFormat = '@'
$ws.Range('G73').Value = 'Pendiente'
$ws.Range('G73').Style = 'Normal'
$ws.Range('H73').NumberFormat = '@'
$ws.Range('H73').Value = 'Columna inclinada'
$ws.Range('H73').Style = 'Normal'
$ws.Range('J73').NumberFormat = '@'
$ws.Range('J73').Value = 'Aplomo'
$ws.Range('J73').Style = 'Normal'
$ws.Range('K73').NumberFormat = '@'
$ws.Range('K73').Value = 'Sin equipos'
$ws.Range('K73').Style = 'Normal'
$ws.Range('L73').NumberFormat = '@'
$ws.Range('L73').Value = 'Terminal'
$ws.Range('L73').Style = 'Normal'
$ws.Range('O73').NumberFormat = '@'
$ws.Range('O73').Value = 'Saavedra'
$ws.Range('O73').Style = 'Normal'
$ws.Range('P73').NumberFormat = '@'
$ws.Range('P73').Value = 'Capital Norte'
$ws.Range('P73').Style = 'Normal'
$ws.Range('I73').Value = 1
$ws.Range('M73').Value = -58.469735
$ws.Range('N73').Value = -34.547232

# Row 74
$ws.Range('A74').NumberFormat = '@'
$ws.Range('A74').Value = '7225'
$ws.Range('A74').Style = 'Normal'
$ws.Range('B74').NumberFormat = '@'
$ws.Range('B74').Value = '9/16/2025'
$ws.Range('B74').Style = 'Normal'
$ws.Range('C74').NumberFormat = '@'
$ws.Range('C74').Value = 'AMENABAR 3590'
$ws.Range('C74').Style = 'Normal'
$ws.Range('D74').NumberFormat = '@'
$ws.Range('D74').Value = '13'
$ws.Range('D74').Style = 'Normal'
$ws.Range('E74').NumberFormat = '@'
$ws.Range('E74').Value = '809784519'
$ws.Range('E74').Style = 'Normal'
$ws.Range('F74').NumberFormat = '@'
$ws.Range('F74').Value = 'PEBCOM'
$ws.Range('F74').Style = 'Normal'
$ws.Range('G74').NumberFormat = '@'
$ws.Range('G74').Value = 'Pendiente'
$ws.Range('G74').Style = 'Normal'
$ws.Range('H74').NumberFormat = '@'
$ws.Range('H74').Value = 'Reparar rienda y tambien reclaman columna picada pero no se ve la foto verificarla y evaluar cambio'
$ws.Range('H74').Style = 'Normal'
$ws.Range('J74').NumberFormat = '@'
$ws.Range('J74').Value = 'Cambio'
$ws.Range('J74').Style = 'Normal'
$ws.Range('K74').NumberFormat = '@'
$ws.Range('K74').Value = 'Sin equipos'
$ws.Range('K74').Style = 'Normal'
$ws.Range('L74').NumberFormat = '@'
$ws.Range('L74').Value = 'Terminal'
$ws.Range('L74').Style = 'Normal'
$ws.Range('O74').NumberFormat = '@'
$ws.Range('O74').Value = 'Saavedra'
$ws.Range('O74').Style = 'Normal'
$ws.Range('P74').NumberFormat = '@'
$ws.Range('P74').Value = 'Capital Norte'
$ws.Range('P74').Style = 'Normal'
$ws.Range('I74').Value = 1
$ws.Range('M74').Value = -58.470045
$ws.Range('N74').Value = -34.550272

# Row 75
$ws.Range('A75').NumberFormat = '@'
$ws.Range('A75').Value = '7234'
$ws.Range('A75').Style = 'Normal'
$ws.Range('B75').NumberFormat = '@'
$ws.Range('B75').Value = '9/16/2025'
$ws.Range('B75').Style = 'Normal'
$ws.Range('C75').NumberFormat = '@'
$ws.Range('C75').Value = 'MOLDES 3388'
$ws.Range('C75').Style = 'Normal'
$ws.Range('D75').NumberFormat = '@'
$ws.Range('D75').Value = '13'
$ws.Range('D75').Style = 'Normal'
$ws.Range('E75').NumberFormat = '@'
$ws.Range('E75').Value = '809784522'
$ws.Range('E75').Style = 'Normal'
$ws.Range('F75').NumberFormat = '@'
$ws.Range('F75').Value = 'PEBCOM'
$ws.Range('F75').Style = 'Normal'
$ws.Range('G75').NumberFormat = '@'
$ws.Range('G75').Value = 'Pendiente'
$ws.Range('G75').Style = 'Normal'
$ws.Range('H75').NumberFormat = '@'
$ws.Range('H75').Value = 'Picada'
$ws.Range('H75').Style = 'Normal'
$ws.Range('J75').NumberFormat = '@'
$ws.Range('J75').Value = 'Cambio'
$ws.Range('J75').Style = 'Normal'
$ws.Range('K75').NumberFormat = '@'
$ws.Range('K75').Value = 'Sin equipos'
$ws.Range('K75').Style = 'Normal'
$ws.Range('L75').NumberFormat = '@'
$ws.Range('L75').Value = 'Pasante'
$ws.Range('L75').Style = 'Normal'
$ws.Range('O75').NumberFormat = '@'
$ws.Range('O75').Value = 'Saavedra'
$ws.Range('O75').Style = 'Normal'
$ws.Range('P75').NumberFormat = '@'
$ws.Range('P75').Value = 'Capital Norte'
$ws.Range('P75').Style = 'Normal'
$ws.Range('I75').Value = 1
$ws.Range('M75').Value = -58.469426
$ws.Range('N75').Value = -34.552639

# Row 76
$ws.Range('A76').NumberFormat = '@'
$ws.Range('A76').Value = '7264'
$ws.Range('A76').Style = 'Normal'
$ws.Range('B76').NumberFormat = '@'
$ws.Range('B76').Value = '9/18/2025'
$ws.Range('B76').Style = 'Normal'
$ws.Range('C76').NumberFormat = '@'
$ws.Range('C76').Value = 'Dorrego 2265'
$ws.Range('C76').Style = 'Normal'
$ws.Range('D76').NumberFormat = '@'
$ws.Range('D76').Value = '14'
$ws.Range('D76').Style = 'Normal'
$ws.Range('E76').NumberFormat = '@'
$ws.Range('E76').Value = '809837500'
$ws.Range('E76').Style = 'Normal'
$ws.Range('F76').NumberFormat = '@'
$ws.Range('F76').Value = 'PEBCOM'
$ws.Range('F76').Style = 'Normal'
$ws.Range('G76').NumberFormat = '@'
$ws.Range('G76').Value = 'Pendiente'
$ws.Range('G76').Style = 'Normal'
$ws.Range('H76').NumberFormat = '@'
$ws.Range('H76').Value = 'Picada'
$ws.Range('H76').Style = 'Normal'
$ws.Range('J76').NumberFormat = '@'
$ws.Range('J76').Value = 'Cambio'
$ws.Range('J76').Style = 'Normal'
$ws.Range('K76').NumberFormat = '@'
$ws.Range('K76').Value = 'Sin equipos'
$ws.Range('K76').Style = 'Normal'
$ws.Range('L76').NumberFormat = '@'
$ws.Range('L76').Value = 'Pasante'
$ws.Range('L76').Style = 'Normal'
$ws.Range('O76').NumberFormat = '@'
$ws.Range('O76').Value = 'Palermo'
$ws.Range('O76').Style = 'Normal'
$ws.Range('P76').NumberFormat = '@'
$ws.Range('P76').Value = 'Capital Sur'
$ws.Range('P76').Style = 'Normal'
$ws.Range('I76').Value = 1
$ws.Range('M76').Value = -58.438083
$ws.Range('N76').Value = -34.577107

# Row 77
$ws.Range('A77').NumberFormat = '@'
$ws.Range('A77').Value = '-603'
$ws.Range('A77').Style = 'Normal'
$ws.Range('B77').NumberFormat = '@'
$ws.Range('B77').Value = '9/22/2025'
$ws.Range('B77').Style = 'Normal'
$ws.Range('C77').NumberFormat = '@'
$ws.Range('C77').Value = 'ANCHORENA, TOMAS MANUEL DE, DR. 821'
$ws.Range('C77').Style = 'Normal'
$ws.Range('D77').NumberFormat = '@'
$ws.Range('D77').Value = '3'
$ws.Range('D77').Style = 'Normal'
$ws.Range('E77').NumberFormat = '@'
$ws.Range('E77').Value = '809910086'
$ws.Range('E77').Style = 'Normal'
$ws.Range('F77').NumberFormat = '@'
$ws.Range('F77').Value = 'PEBCOM'
$ws.Range('F77').Style = 'Normal'
$ws.Range('G77').NumberFormat = '@'
$ws.Range('G77').Value = 'Pendiente'
$ws.Range('G77').Style = 'Normal'
$ws.Range('H77').NumberFormat = '@'
$ws.Range('H77').Value = 'Columna chocada pendiente para instalar un corporativo'
$ws.Range('H77').Style = 'Normal'
$ws.Range('J77').NumberFormat = '@'
$ws.Range('J77').Value = 'Cambio'
$ws.Range('J77').Style = 'Normal'
$ws.Range('K77').NumberFormat = '@'
$ws.Range('K77').Value = 'Sin equipos'
$ws.Range('K77').Style = 'Normal'
$ws.Range('L77').NumberFormat = '@'
$ws.Range('L77').Value = 'Pasante'
$ws.Range('L77').Style = 'Normal'
$ws.Range('O77').NumberFormat = '@'
$ws.Range('O77').Value = 'Almagro'
$ws.Range('O77').Style = 'Normal'
$ws.Range('P77').NumberFormat = '@'
$ws.Range('P77').Value = 'Capital Sur'
$ws.Range('P77').Style = 'Normal'
$ws.Range('I77').Value = 1
$ws.Range('M77').Value = -58.408551
$ws.Range('N77').Value = -34.599265

# Row 78
$ws.Range('A78').NumberFormat = '@'
$ws.Range('A78').Value = '-612'
$ws.Range('A78').Style = 'Normal'
$ws.Range('B78').NumberFormat = '@'
$ws.Range('B78').Value = '9/24/2025'
$ws.Range('B78').Style = 'Normal'
$ws.Range('C78').NumberFormat = '@'
$ws.Range('C78').Value = 'Herrera 588'
$ws.Range('C78').Style = 'Normal'
$ws.Range('D78').NumberFormat = '@'
$ws.Range('D78').Value = '4'
$ws.Range('D78').Style = 'Normal'
$ws.Range('E78').NumberFormat = '@'
$ws.Range('E78').Value = '809972811'
$ws.Range('E78').Style = 'Normal'
$ws.Range('F78').NumberFormat = '@'
$ws.Range('F78').Value = 'PEBCOM'
$ws.Range('F78').Style = 'Normal'
$ws.Range('G78').NumberFormat = '@'
$ws.Range('G78').Value = 'Pendiente'
$ws.Range('G78').Style = 'Normal'
$ws.Range('H78').NumberFormat = '@'
$ws.Range('H78').Value = ''
$ws.Range('H78').Style = 'Normal'
$ws.Range('J78').NumberFormat = '@'
$ws.Range('J78').Value = 'Cambio'
$ws.Range('J78').Style = 'Normal'
$ws.Range('K78').NumberFormat = '@'
$ws.Range('K78').Value = 'Sin equipos'
$ws.Range('K78').Style = 'Normal'
$ws.Range('L78').NumberFormat = '@'
$ws.Range('L78').Value = 'Pasante'
$ws.Range('L78').Style = 'Normal'
$ws.Range('O78').NumberFormat = '@'
$ws.Range('O78').Value = 'San Telmo'
$ws.Range('O78').Style = 'Normal'
$ws.Range('P78').NumberFormat = '@'
$ws.Range('P78').Value = 'Capital Sur'
$ws.Range('P78').Style = 'Normal'
$ws.Range('I78').Value = 1
$ws.Range('M78').Value = -58.378275
$ws.Range('N78').Value = -34.635935

# Row 79
$ws.Range('A79').NumberFormat = '@'
$ws.Range('A79').Value = '-613'
$ws.Range('A79').Style = 'Normal'
$ws.Range('B79').NumberFormat = '@'
$ws.Range('B79').Value = '9/24/2025'
$ws.Range('B79').Style = 'Normal'
$ws.Range('C79').NumberFormat = '@'
$ws.Range('C79').Value = 'Aristobulo del valle 1875'
$ws.Range('C79').Style = 'Normal'
$ws.Range('D79').NumberFormat = '@'
$ws.Range('D79').Value = '4'
$ws.Range('D79').Style = 'Normal'
$ws.Range('E79').NumberFormat = '@'
$ws.Range('E79').Value = '809972816'
$ws.Range('E79').Style = 'Normal'
$ws.Range('F79').NumberFormat = '@'
$ws.Range('F79').Value = 'PEBCOM'
$ws.Range('F79').Style = 'Normal'
$ws.Range('G79').NumberFormat = '@'
$ws.Range('G79').Value = 'Pendiente'
$ws.Range('G79').Style = 'Normal'
$ws.Range('H79').NumberFormat = '@'
$ws.Range('H79').Value = ''
$ws.Range('H79').Style = 'Normal'
$ws.Range('J79').NumberFormat = '@'
$ws.Range('J79').Value = 'Cambio'
$ws.Range('J79').Style = 'Normal'
$ws.Range('K79').NumberFormat = '@'
$ws.Range('K79').Value = 'Sin equipos'
$ws.Range('K79').Style = 'Normal'
$ws.Range('L79').NumberFormat = '@'
$ws.Range('L79').Value = 'Pasante'
$ws.Range('L79').Style = 'Normal'
$ws.Range('O79').NumberFormat = '@'
$ws.Range('O79').Value = 'San Telmo'
$ws.Range('O79').Style = 'Normal'
$ws.Range('P79').NumberFormat = '@'
$ws.Range('P79').Value = 'Capital Sur'
$ws.Range('P79').Style = 'Normal'
$ws.Range('I79').Value = 0
$ws.Range('M79').Value = -58.377585
$ws.Range('N79').Value = -34.636595
